$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use an existing unstyled (default style) cell as a style reference so that
# assigning text values to these numeric-looking cells does not introduce a
# new explicit cell style (they must stay on the default/general style).
$refStyle = $ws.Range("F2").Style

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.Value = "'" + $value
    $rng.Style = $refStyle
}

Set-TextValue "D2" "289.54"
Set-TextValue "E2" "-0.31%"
Set-TextValue "D3" "31.07"
Set-TextValue "E3" "1.53%"
Set-TextValue "D4" "4.952"
Set-TextValue "E4" "0.76%"
Set-TextValue "D5" "0.07360"
Set-TextValue "D6" "2.310"
Set-TextValue "E6" "27.60%"
Set-TextValue "D7" "7.678"
Set-TextValue "E7" "0.50%"
Set-TextValue "D8" "0.9194"
Set-TextValue "E8" "2.12%"
Set-TextValue "D9" "0.09111"
Set-TextValue "E9" "13.56%"
Set-TextValue "D10" "0.1706"
Set-TextValue "E10" "1.44%"
Set-TextValue "D11" "0.08184"
Set-TextValue "E11" "1.61%"
Set-TextValue "D12" "0.03113"
Set-TextValue "E12" "1.71%"
Set-TextValue "D13" "0.09994"
Set-TextValue "E13" "-0.10%"
Set-TextValue "D14" "0.001499"
Set-TextValue "E14" "-0.40%"
Set-TextValue "D15" "0.005733"
Set-TextValue "E15" "-0.16%"
Set-TextValue "E16" "-0.13%"
Set-TextValue "D17" "3.744"
Set-TextValue "E17" "1.28%"
Set-TextValue "D18" "2.029"
Set-TextValue "E18" "-2.28%"
Set-TextValue "E19" "0.42%"
Set-TextValue "D20" "0.1299"
Set-TextValue "E20" "-0.32%"
Set-TextValue "D21" "4.181"
Set-TextValue "E21" "5.70%"
Set-TextValue "D22" "0.2124"
Set-TextValue "E22" "-1.97%"
Set-TextValue "D23" "0.04511"
Set-TextValue "E23" "0.13%"
Set-TextValue "E24" "0.15%"
Set-TextValue "D25" "0.004199"
Set-TextValue "E25" "-5.28%"
Set-TextValue "E26" "0.21%"
Set-TextValue "D27" "0.0003394"
Set-TextValue "E27" "0.18%"
Set-TextValue "D39" "0.01579"
Set-TextValue "E39" "-0.81%"
Set-TextValue "D40" "0.04514"
Set-TextValue "E40" "3.82%"
Set-TextValue "D41" "0.007376"
Set-TextValue "E41" "0.74%"
Set-TextValue "D42" "0.009851"
Set-TextValue "E42" "-1.72%"
Set-TextValue "D43" "0.1338"
Set-TextValue "E43" "1.94%"
Set-TextValue "E44" "10.96%"
Set-TextValue "D45" "0.008503"
Set-TextValue "E45" "-10.10%"
Set-TextValue "D46" "0.00006104"
Set-TextValue "E46" "3.96%"
Set-TextValue "E47" "0.20%"
Set-TextValue "D48" "2.606"
Set-TextValue "E48" "15.58%"
Set-TextValue "D49" "0.002000"
Set-TextValue "E49" "-30.88%"
Set-TextValue "D50" "0.00002100"
Set-TextValue "E50" "0.20%"
Set-TextValue "D51" "0.0002000"
Set-TextValue "E51" "0.20%"
